# Updates the crypto price/volume table to the "16-12-2022 08" snapshot.
# Price column (D) cells are stored as text (they hold values like "--" for
# some rows), so numeric-looking updates are written with a leading
# apostrophe to force text entry instead of letting Excel coerce them to
# numbers. Coin/Link/Volume columns (B, C, E) are plain text already.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'262.82"
$ws.Range("D3").Value = "'24.17"
$ws.Range("D4").Value = "'6.200"
$ws.Range("D5").Value = "'0.06232"
$ws.Range("D6").Value = "'6.743"
$ws.Range("D7").Value = "'3.444"
$ws.Range("D8").Value = "'1.353"
$ws.Range("D9").Value = "'0.8011"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1593"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08129"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03396"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03092"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09344"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.703"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001712"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04810"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006152"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006224"
$ws.Range("D20").Value = "'0.006178"
$ws.Range("D21").Value = "'0.001100"
$ws.Range("D22").Value = "'0.0001500"
$ws.Range("D23").Value = "'3.688"
$ws.Range("D25").Value = "'0.3375"
$ws.Range("D26").Value = "'0.1275"
$ws.Range("D27").Value = "'0.0003203"
$ws.Range("D40").Value = "'0.04637"
$ws.Range("D41").Value = "'0.007051"
$ws.Range("D43").Value = "'0.003601"
$ws.Range("D44").Value = "'0.01000"
$ws.Range("D45").Value = "'0.002971"
$ws.Range("D46").Value = "'0.00005877"
$ws.Range("D48").Value = "'0.7003"
$ws.Range("D49").Value = "'0.1283"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D51").Value = "'0.01010"
